# Autogenerated on Wed Apr 01 2015 00:15:40 GMT+0000 (Coordinated Universal Time)
#
# Refresh the Malta MSME "Summary" sheet with updated statistics
# (Enterprises density, Employment % of total, Enterprises % of total).
# The underlying cells store these figures as text, so each cell is
# briefly switched to a text number format while the new value is
# entered (otherwise Excel would auto-convert the numeric-looking
# string into a real number) and then restored to the "Normal" style
# so no visible formatting changes remain.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B13" = "59.37"   # Enterprises density (per 1000 people) - Micro
    "C13" = "3.61"    # Enterprises density (per 1000 people) - SMEs
    "D13" = "62.97"   # Enterprises density (per 1000 people) - MSMEs

    "B14" = "35.97"   # Employment (% of total) - Micro
    "C14" = "44.07"   # Employment (% of total) - SMEs
    "D14" = "80.03"   # Employment (% of total) - MSMEs

    "B16" = "94.11"   # Enterprises (% of total) - Micro
    "C16" = "5.72"    # Enterprises (% of total) - SMEs
    "D16" = "99.83"   # Enterprises (% of total) - MSMEs
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
